$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247; everything from 247..315 shifts down to 248..316.
$ws.Rows(247).Insert()

# Fill the newly inserted row 247 with the new weekly price-report entry.
$ws.Range("A247").Value = 4
$ws.Range("B247").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C247").Value = "Los Lagos"
$ws.Range("D247").Value = 44841
$ws.Range("E247").Value = 10
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100108
$ws.Range("H247").Value = "Tropicales y subtropicales"
$ws.Range("I247").Value = 100108005
$ws.Range("J247").Value = "Piña"
$ws.Range("K247").Value = "Caramelo"
$ws.Range("L247").Value = "Primera"
$ws.Range("M247").Value = 200
$ws.Range("N247").Value = 23000
$ws.Range("O247").Value = 23500
$ws.Range("P247").Value = 23250
$ws.Range("Q247").Value = "$/caja 12 unidades"
$ws.Range("R247").Value = "Ecuador"
$ws.Range("S247").Value = 1938
$ws.Range("T247").Value = 12
